$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 22727760
$ws.Range("I12").Value = 635.1429000000001
$ws.Range("K12").Value = 635.1429000000001
$ws.Range("M12").Value = -465.1429000000001
$ws.Range("H17").Value = 2223.625
$ws.Range("J17").Value = 2223.625
$ws.Range("L17").Value = 6670.875
$ws.Range("N17").Value = -7006.875
$ws.Range("H40").Value = 11640.647
$ws.Range("J40").Value = 11640.647
$ws.Range("L40").Value = 11640.647
$ws.Range("N40").Value = -11990.647
$ws.Range("H64").Value = 4944.4614
$ws.Range("I64").Value = 4879.6665
$ws.Range("K64").Value = 4879.6665
$ws.Range("M64").Value = -4631.6665
$ws.Range("H67").Value = 4944.4614
$ws.Range("I67").Value = 4879.6665
$ws.Range("K67").Value = 4879.6665
$ws.Range("M67").Value = -4021.6665
$ws.Range("H86").Value = 3874.111
$ws.Range("I86").Value = 1594.9231
$ws.Range("J86").Value = 9800
$ws.Range("K86").Value = 1594.9231
$ws.Range("L86").Value = 9800
$ws.Range("M86").Value = -471.9231
$ws.Range("N86").Value = -12046
$ws.Range("H89").Value = 3874.111
$ws.Range("I89").Value = 1594.9231
$ws.Range("J89").Value = 9800
$ws.Range("K89").Value = 7974.6155
$ws.Range("L89").Value = 49000
$ws.Range("M89").Value = -2358.6155
$ws.Range("N89").Value = -60232
$ws.Range("H101").Value = 1310.75
$ws.Range("I101").Value = 1310.75
$ws.Range("K101").Value = 3932.25
$ws.Range("M101").Value = -2310.25
$ws.Range("H127").Value = 2147.5
$ws.Range("I127").Value = 2147.5
$ws.Range("K127").Value = 6442.5
$ws.Range("M127").Value = -1482.5
$ws.Range("H138").Value = 2277.2666
$ws.Range("I138").Value = 999.1
$ws.Range("J138").Value = 3299.8
$ws.Range("K138").Value = 2997.3
$ws.Range("L138").Value = 9899.400000000001
$ws.Range("M138").Value = 2142.7
$ws.Range("N138").Value = -20179.4
$ws.Range("H139").Value = 90825
$ws.Range("J139").Value = 90825
$ws.Range("L139").Value = 90825
$ws.Range("N139").Value = -101105

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7768.044
$ws.Range("I32").Value = 3102.574
$ws.Range("K32").Value = 3102.574
$ws.Range("M32").Value = -2815.574
$ws.Range("H102").Value = 69251.586
$ws.Range("I102").Value = 73144
$ws.Range("K102").Value = 73144
$ws.Range("M102").Value = -71522
$ws.Range("H132").Value = 2396.625
$ws.Range("I132").Value = 1814.3125
$ws.Range("J132").Value = 3561.25
$ws.Range("K132").Value = 5442.9375
$ws.Range("L132").Value = 10683.75
$ws.Range("M132").Value = -2912.9375
$ws.Range("N132").Value = -15743.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 196
$ws.Range("J19").Value = 196
$ws.Range("L19").Value = 196
$ws.Range("N19").Value = -542
$ws.Range("H62").Value = 75000
$ws.Range("I62").Value = 75000
$ws.Range("K62").Value = 75000
$ws.Range("M62").Value = -74314
$ws.Range("H65").Value = 75000
$ws.Range("I65").Value = 75000
$ws.Range("K65").Value = 225000
$ws.Range("M65").Value = -221568
$ws.Range("H86").Value = 3928.1667
$ws.Range("I86").Value = 3114
$ws.Range("K86").Value = 3114
$ws.Range("M86").Value = -1991
$ws.Range("H89").Value = 3928.1667
$ws.Range("I89").Value = 3114
$ws.Range("K89").Value = 15570
$ws.Range("M89").Value = -9954
$ws.Range("H99").Value = 1698939.8
$ws.Range("I99").Value = 64112.938
$ws.Range("J99").Value = 10418016
$ws.Range("K99").Value = 64112.938
$ws.Range("L99").Value = 10418016
$ws.Range("M99").Value = -62614.938
$ws.Range("N99").Value = -10421012
$ws.Range("H105").Value = 57918.332
$ws.Range("I105").Value = 78386.16
$ws.Range("J105").Value = 4702
$ws.Range("K105").Value = 78386.16
$ws.Range("L105").Value = 4702
$ws.Range("M105").Value = -76639.16
$ws.Range("N105").Value = -8196
$ws.Range("H107").Value = 3735.7144
$ws.Range("I107").Value = 2310.375
$ws.Range("K107").Value = 2310.375
$ws.Range("M107").Value = -390.375
$ws.Range("H132").Value = 27017.596
$ws.Range("J132").Value = 27017.596
$ws.Range("L132").Value = 27017.596
$ws.Range("N132").Value = -37137.59600000001
$ws.Range("H138").Value = 79950.664
$ws.Range("J138").Value = 79950.664
$ws.Range("L138").Value = 79950.664
$ws.Range("N138").Value = -90230.664

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 2999
$ws.Range("J29").Value = 2999
$ws.Range("L29").Value = 2999
$ws.Range("N29").Value = -3585
$ws.Range("H31").Value = 12520.108
$ws.Range("I31").Value = 2648.4211
$ws.Range("K31").Value = 2648.4211
$ws.Range("M31").Value = -2353.4211
$ws.Range("H34").Value = 12520.108
$ws.Range("I34").Value = 2648.4211
$ws.Range("K34").Value = 2648.4211
$ws.Range("M34").Value = -2446.4211
$ws.Range("H94").Value = 1026.24
$ws.Range("I94").Value = 1017.1
$ws.Range("J94").Value = 1032.3334
$ws.Range("K94").Value = 1017.1
$ws.Range("L94").Value = 1032.3334
$ws.Range("M94").Value = -566.1
$ws.Range("N94").Value = -1934.3334
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 2937.25
$ws.Range("J46").Value = 924.5
$ws.Range("L46").Value = 2773.5
$ws.Range("N46").Value = -2955.5
$ws.Range("H132").Value = 2708.75
$ws.Range("I132").Value = 2361
$ws.Range("J132").Value = 2957.1428
$ws.Range("K132").Value = 21249
$ws.Range("L132").Value = 26614.2852
$ws.Range("M132").Value = -18719
$ws.Range("N132").Value = -31674.2852

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 127080.21
$ws.Range("I70").Value = 149574.14
$ws.Range("J70").Value = 104586.29
$ws.Range("K70").Value = 149574.14
$ws.Range("L70").Value = 104586.29
$ws.Range("M70").Value = -149304.14
$ws.Range("N70").Value = -105126.29
$ws.Range("H73").Value = 127080.21
$ws.Range("I73").Value = 149574.14
$ws.Range("J73").Value = 104586.29
$ws.Range("K73").Value = 149574.14
$ws.Range("L73").Value = 104586.29
$ws.Range("M73").Value = -148638.14
$ws.Range("N73").Value = -106458.29
$ws.Range("H136").Value = 38869.2
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 38869.2
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 116607.6
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -121707.6
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1063.7273
$ws.Range("I22").Value = 821.2857
$ws.Range("J22").Value = 1488
$ws.Range("K22").Value = 821.2857
$ws.Range("L22").Value = 1488
$ws.Range("M22").Value = -526.2857
$ws.Range("N22").Value = -2078
$ws.Range("H27").Value = 1063.7273
$ws.Range("I27").Value = 821.2857
$ws.Range("J27").Value = 1488
$ws.Range("K27").Value = 821.2857
$ws.Range("L27").Value = 1488
$ws.Range("M27").Value = -714.2857
$ws.Range("N27").Value = -1702
$ws.Range("H40").Value = 5557851
$ws.Range("I40").Value = 2213
$ws.Range("K40").Value = 2213
$ws.Range("M40").Value = -2077
$ws.Range("H46").Value = 2654.25
$ws.Range("I46").Value = 1598.5
$ws.Range("K46").Value = 1598.5
$ws.Range("M46").Value = -1410.5
$ws.Range("H55").Value = 1870.1936
$ws.Range("I55").Value = 690.8
$ws.Range("J55").Value = 2975.875
$ws.Range("K55").Value = 690.8
$ws.Range("L55").Value = 2975.875
$ws.Range("M55").Value = -517.8
$ws.Range("N55").Value = -3321.875
$ws.Range("H93").Value = 914.7778
$ws.Range("I93").Value = 965.9
$ws.Range("J93").Value = 850.875
$ws.Range("K93").Value = 965.9
$ws.Range("L93").Value = 850.875
$ws.Range("M93").Value = 282.1
$ws.Range("N93").Value = -3346.875
$ws.Range("H122").Value = 11793717
$ws.Range("I122").Value = 36706.31
$ws.Range("J122").Value = 50004000
$ws.Range("K122").Value = 110118.93
$ws.Range("L122").Value = 150012000
$ws.Range("M122").Value = -107668.93
$ws.Range("N122").Value = -150016900

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 145328.89
$ws.Range("J46").Value = 145328.89
$ws.Range("L46").Value = 145328.89
$ws.Range("N46").Value = -145790.89
$ws.Range("H81").Value = 3564.5
$ws.Range("I81").Value = 2502.9092
$ws.Range("K81").Value = 5005.8184
$ws.Range("M81").Value = -3944.8184
$ws.Range("H84").Value = 3564.5
$ws.Range("I84").Value = 2502.9092
$ws.Range("K84").Value = 25029.092
$ws.Range("M84").Value = -19725.092
$ws.Range("H100").Value = 7145185
$ws.Range("J100").Value = 1775
$ws.Range("L100").Value = 3550
$ws.Range("N100").Value = -4632
$ws.Range("H107").Value = 7885.6045
$ws.Range("I107").Value = 7568.4287
$ws.Range("K107").Value = 22705.2861
$ws.Range("M107").Value = -20785.2861
$ws.Range("H113").Value = 669.1
$ws.Range("I113").Value = 490.35715
$ws.Range("K113").Value = 1471.07145
$ws.Range("M113").Value = 698.9285500000001
$ws.Range("H132").Value = 1673803
$ws.Range("I132").Value = 1190.5625
$ws.Range("K132").Value = 3571.6875
$ws.Range("M132").Value = -1041.6875
$ws.Range("H134").Value = 145328.89
$ws.Range("J134").Value = 145328.89
$ws.Range("L134").Value = 435986.67
$ws.Range("N134").Value = -441056.67
